$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new "Material" column before the old Class column (C), shifting
# Class/Progress/the summary formulas one column to the right.
$ws.Columns.Item(3).Insert()

# New column header + per-chapter material names (chapters 1-6, rows 2-7).
$ws.Range("C1").Value = "Material"
$ws.Range("C2").Value = "Binary"
$ws.Range("C3").Value = "Array"
$ws.Range("C4").Value = "String"
$ws.Range("C5").Value = "Linked L"
$ws.Range("C6").Value = "Stacks"
$ws.Range("C7").Value = "Binary T"

# Chapter 5 (row 6) progress: 5 more questions finished.
$ws.Range("E5").Value = 13
$ws.Range("E6").Value = 4

# Drop the stray "half of total problems" helper formula that used to live
# at G6 (now shifted to H6 after the column insert).
$ws.Range("H6").ClearContents()

# Fix the chapter-number column: it had a duplicate "4" at row 6 (old chapter
# numbering bug). Replace the static numbers from row 6 down with an
# auto-incrementing formula instead.
$ws.Range("A6").Formula = "=A5+1"
$ws.Range("A7:A22").Formula = "=A6+1"

# Selection moves to D24 after the edit.
$ws.Range("D24").Select()
